$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 281 - this shifts the former rows 281-301
# down to 282-302 and grows the sheet's used range by one row.
$ws.Rows.Item(281).Insert()

# Populate the newly inserted row 281 with the new weekly record
# (Terminal Hortofrutícola Agro Chillán - Piña, Caramelo, Primera).
$ws.Cells.Item(281, 1).Value = 7
$ws.Cells.Item(281, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(281, 3).Value = "Ñuble"
$ws.Cells.Item(281, 4).Value = 45021
$ws.Cells.Item(281, 5).Value = 16
$ws.Cells.Item(281, 6).Value = "Fruta"
$ws.Cells.Item(281, 7).Value = 100108
$ws.Cells.Item(281, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(281, 9).Value = 100108005
$ws.Cells.Item(281, 10).Value = "Piña"
$ws.Cells.Item(281, 11).Value = "Caramelo"
$ws.Cells.Item(281, 12).Value = "Primera"
$ws.Cells.Item(281, 13).Value = 60
$ws.Cells.Item(281, 14).Value = 20000
$ws.Cells.Item(281, 15).Value = 20000
$ws.Cells.Item(281, 16).Value = 20000
$ws.Cells.Item(281, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(281, 18).Value = "Ecuador"
$ws.Cells.Item(281, 19).Value = 1667
$ws.Cells.Item(281, 20).Value = 12
